$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Texas -- Bexar County" -- date published moves forward one day ---
$ws.Range("B3").Value = 44034

# --- Row 4: "New York -- New York" -- this run's scrape failed for NY, so the
#     previously-successful values are wiped and the status cell now carries
#     the rate-limit error instead of "Success!" ---
$ws.Range("B4:H4").ClearContents()
$ws.Range("B4:H4").ClearFormats()
$ws.Range("K4:L4").ClearContents()
$ws.Range("K4:L4").ClearFormats()

$ws.Range("I4").Value = $false
$ws.Range("J4").Value = $false

$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"

# --- Row 39: "Delaware" -- different scrape failure this run ---
$ws.Range("O39").Value = "An error occurred. ... HTTPError('504 Server Error: Gateway Time-out for url: https://myhealthycommunity.dhss.delaware.gov/locations/state/')"
